$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Red Status count
$ws.Range("G1").Value = "Red Status: 14 projects"

# Add new data points in column G for rows 8-15
$ws.Range("G8").Value  = "q1: 23"
$ws.Range("G9").Value  = "q2: 24"
$ws.Range("G10").Value = "q3: 25"
$ws.Range("G11").Value = "q4: 26"
$ws.Range("G12").Value = "q5: 27"
$ws.Range("G13").Value = "q6: 28"
$ws.Range("G14").Value = "q7: 29"
$ws.Range("G15").Value = "q8: 30"
